# Added 4wk low sales check
# Update "Forecast Comparison" sheet with recalculated forecast metrics
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# Row -> D (MyForecast), H (Inventory Coverage), I (Stockout Risk), J (Reorder Urgency), L (Seasonality Index)
$rows = @(
    @{ Row = 2;  D = 0; H = $null; I = "Low"; J = "Normal"; L = 1.18 },
    @{ Row = 3;  D = 0; H = $null; I = "Low"; J = "Normal"; L = 1.2 },
    @{ Row = 4;  D = 0; H = $null; I = "Low"; J = "Normal"; L = 1.13 },
    @{ Row = 5;  D = 0; H = $null; I = "Low"; J = "Normal"; L = 1.05 },
    @{ Row = 6;  D = 0; H = $null; I = "Low"; J = "Normal"; L = 0.85 },
    @{ Row = 7;  D = 0; H = $null; I = "Low"; J = "Normal"; L = 0.96 },
    @{ Row = 8;  D = 0; H = $null; I = "Low"; J = "Normal"; L = 1.06 },
    @{ Row = 9;  D = 0; H = $null; I = "Low"; J = "Normal"; L = 0.9399999999999999 },
    @{ Row = 10; D = 0; H = $null; I = "Low"; J = "Normal"; L = 0.88 },
    @{ Row = 11; D = 0; H = $null; I = "Low"; J = "Normal"; L = 0.92 },
    @{ Row = 12; D = 0; H = $null; I = "Low"; J = "Normal"; L = 1.11 },
    @{ Row = 13; D = 0; H = $null; I = "Low"; J = "Normal"; L = 1.13 },
    @{ Row = 14; D = 0; H = $null; I = "Low"; J = "Normal"; L = 0.9 },
    @{ Row = 15; D = 0; H = $null; I = "Low"; J = "Normal"; L = 1.14 },
    @{ Row = 16; D = 0; H = $null; I = "Low"; J = "Normal"; L = 0.87 },
    @{ Row = 17; D = 1; H = 7;     I = "Low"; J = "Normal"; L = 0.83 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws1.Cells.Item($n, 4).Value = $r.D
    if ($null -eq $r.H) {
        $ws1.Cells.Item($n, 8).Value = ""
    } else {
        $ws1.Cells.Item($n, 8).Value = $r.H
    }
    $ws1.Cells.Item($n, 9).Value = $r.I
    $ws1.Cells.Item($n, 10).Value = $r.J
    $ws1.Cells.Item($n, 12).Value = $r.L
}

# Update "Summary" sheet totals (kept as text to match existing column formatting)
$summaryUpdates = @(
    @{ Cell = "B9";  Value = "1" },
    @{ Cell = "B10"; Value = "0" },
    @{ Cell = "B11"; Value = "0" },
    @{ Cell = "B12"; Value = "1" },
    @{ Cell = "B14"; Value = "0" }
)
foreach ($u in $summaryUpdates) {
    $cell = $ws2.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
